$d = $word.ActiveDocument

# Paragraph 56: "Micro desafío - code review" heading.
# Paragraph 58: "Repasando todo lo aprendido en Flexbox, ... (quedaron uno
#               al lado del otro pero no de la forma que esta en la
#               imagen) Además, cada vez que se pase el mouse ..."
$headingPara = $d.Paragraphs(56)
$bodyPara = $d.Paragraphs(58)

# --- 1. Text restructuring: "pero no " + "de" + " la forma que " runs merge ---
# (visible text is unchanged: " pero no de la forma que " stays the same,
# but re-running it through Find/Replace re-splits the underlying runs
# exactly like the target revision: the "de" run disappears into the
# previous run, and the " la forma que " run disappears, leaving "esta"
# (the ex-spellchecked run) as the following run.) Scoped to the body
# paragraph so it cannot touch the unrelated "pero no funciona)" text
# elsewhere in the document.
$bodyPara.Range.Find.Execute("pero no de la forma que", $false, $false, $false, $false, $false, $true, 1, $false, "pero no de la forma que", 2) | Out-Null

# --- 2. Apply strikethrough ---
# 2a. The heading paragraph (including its paragraph mark) gets struck
#     through in full.
$headingPara.Range.Font.StrikeThrough = 1

# 2b. The body paragraph gets struck through only for the part preceding
#     "Además" (the rest already carries strike formatting), up to and
#     including the closing parenthesis of "...en la imagen)" - the
#     single space right before "Además" stays unstruck.
$searchRange = $bodyPara.Range.Duplicate
$searchRange.Find.Execute("en la imagen)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target = $d.Range($bodyPara.Range.Start, $searchRange.End)
$target.Font.StrikeThrough = 1
